# Update crypto price/volume table cells per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.612.93"
$ws.Range("E2").Value = "'  +2.39%  "
$ws.Range("D3").Value = "'1.859.89"
$ws.Range("E3").Value = "'  +1.51%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'244.94"
$ws.Range("E5").Value = "'  +1.66%  "
$ws.Range("D6").Value = "'0.6954"
$ws.Range("E6").Value = "'  +1.07%  "
$ws.Range("D7").Value = "'0.9998"
$ws.Range("D8").Value = "'0.07705"
$ws.Range("E8").Value = "'  +0.63%  "
$ws.Range("E9").Value = "'  +0.34%  "
$ws.Range("D10").Value = "'23.75"
$ws.Range("E10").Value = "'  +0.73%  "
$ws.Range("D11").Value = "'0.07764"
$ws.Range("E11").Value = "'  -0.62%  "
$ws.Range("D12").Value = "'5.157"
$ws.Range("E12").Value = "'  +1.71%  "
$ws.Range("D13").Value = "'1.863.39"
$ws.Range("E13").Value = "'  +1.77%  "
$ws.Range("D14").Value = "'92.01"
$ws.Range("E14").Value = "'  +1.76%  "
$ws.Range("D15").Value = "'0.6923"
$ws.Range("E15").Value = "'  +2.29%  "
$ws.Range("D16").Value = "'6.574"
$ws.Range("E16").Value = "'  +2.06%  "
$ws.Range("D17").Value = "'29.603.43"
$ws.Range("E17").Value = "'  +2.39%  "
$ws.Range("D18").Value = "'0.000008311"
$ws.Range("E18").Value = "'  +0.50%  "
$ws.Range("D19").Value = "'2.101.62"
$ws.Range("E19").Value = "'  +1.38%  "
$ws.Range("D20").Value = "'240.81"
$ws.Range("E20").Value = "'  -0.83%  "
$ws.Range("E21").Value = "'  +0.81%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("E22").Value = "'  +0.02%  "
$ws.Range("D23").Value = "'7.603"
$ws.Range("E23").Value = "'  +2.43%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "'  +0.09%  "
$ws.Range("E25").Value = "'  +1.90%  "
$ws.Range("D26").Value = "'8.927"
$ws.Range("E26").Value = "'  +1.69%  "
$ws.Range("D27").Value = "'160.13"
$ws.Range("E27").Value = "'  -0.76%  "
$ws.Range("D28").Value = "'18.29"
$ws.Range("E28").Value = "'  +0.52%  "
$ws.Range("D29").Value = "'1.537"
$ws.Range("E29").Value = "'  +0.06%  "
$ws.Range("D30").Value = "'4.251"
$ws.Range("E30").Value = "'  +0.86%  "
$ws.Range("D31").Value = "'4.187"
$ws.Range("E31").Value = "'  +1.71%  "
$ws.Range("E32").Value = "'  +0.22%  "
$ws.Range("E33").Value = "'  -0.51%  "
$ws.Range("D34").Value = "'0.7735"
$ws.Range("E34").Value = "'  +2.54%  "
$ws.Range("D35").Value = "'1.893"
$ws.Range("E35").Value = "'  +3.45%  "
$ws.Range("D36").Value = "'1.153"
$ws.Range("E36").Value = "'  +0.71%  "
$ws.Range("D37").Value = "'2.686"
$ws.Range("E37").Value = "'  +0.50%  "
$ws.Range("D38").Value = "'1.330.43"
$ws.Range("E38").Value = "'  +8.22%  "
$ws.Range("E39").Value = "'  +1.40%  "
$ws.Range("D40").Value = "'2.727"
$ws.Range("E40").Value = "'  +1.40%  "
$ws.Range("D41").Value = "'0.9673"
$ws.Range("E41").Value = "'  +5.27%  "
$ws.Range("D42").Value = "'106.81"
$ws.Range("E42").Value = "'  -1.37%  "
$ws.Range("D43").Value = "'5.807"
$ws.Range("E43").Value = "'  +4.52%  "
$ws.Range("D44").Value = "'0.9996"
$ws.Range("E44").Value = "'  +0.14%  "
$ws.Range("B45").Value = "'EnergySwap"
$ws.Range("C45").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'9.774"
$ws.Range("E45").Value = "'  +2.79%  "
$ws.Range("B46").Value = "'BabyDogeCoin"
$ws.Range("C46").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000125"
$ws.Range("E46").Value = "'  +3.01%  "
$ws.Range("D47").Value = "'2.000.92"
$ws.Range("E47").Value = "'  +1.31%  "
$ws.Range("E48").Value = "'  +0.86%  "
$ws.Range("D49").Value = "'1.777"
$ws.Range("E49").Value = "'  +2.33%  "
$ws.Range("D50").Value = "'63.61"
$ws.Range("E50").Value = "'  -0.58%  "
$ws.Range("D51").Value = "'6.959"
$ws.Range("E51").Value = "'  +1.02%  "
